$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = "GK_x"
$ws.Cells.Item(2, 2).Value = 0.2123405037926963
$ws.Cells.Item(2, 3).Value = -16.59701423989972
$ws.Cells.Item(2, 4).Value = 0.2123405037926963
$ws.Cells.Item(3, 1).Value = "tpr_x"
$ws.Cells.Item(3, 2).Value = 0.01272532166281927
$ws.Cells.Item(3, 3).Value = -16.59701423989972
$ws.Cells.Item(3, 4).Value = 0.01272532166281927
$ws.Cells.Item(4, 1).Value = "Anticipation_x"
$ws.Cells.Item(4, 2).Value = 0.012991908494656
$ws.Cells.Item(4, 3).Value = -16.59701423989972
$ws.Cells.Item(4, 4).Value = 0.012991908494656
$ws.Cells.Item(5, 1).Value = "Teamwork_x"
$ws.Cells.Item(5, 2).Value = 0.06608506427750496
$ws.Cells.Item(5, 3).Value = -16.59701423989972
$ws.Cells.Item(5, 4).Value = 0.06608506427750496
$ws.Cells.Item(6, 1).Value = "Corners_x"
$ws.Cells.Item(6, 2).Value = 0.04236803143712155
$ws.Cells.Item(6, 3).Value = -16.59701423989972
$ws.Cells.Item(6, 4).Value = 0.04236803143712155
$ws.Cells.Item(7, 1).Value = "Dribbling_x"
$ws.Cells.Item(7, 2).Value = -0.05819302642662285
$ws.Cells.Item(7, 3).Value = -16.59701423989972
$ws.Cells.Item(7, 4).Value = 0.05819302642662285
$ws.Cells.Item(8, 1).Value = "Agility_x"
$ws.Cells.Item(8, 2).Value = 0.01822184493229715
$ws.Cells.Item(8, 3).Value = -16.59701423989972
$ws.Cells.Item(8, 4).Value = 0.01822184493229715
$ws.Cells.Item(9, 1).Value = "Off_the_Ball_x"
$ws.Cells.Item(9, 2).Value = 0.09703207705940729
$ws.Cells.Item(9, 3).Value = -16.59701423989972
$ws.Cells.Item(9, 4).Value = 0.09703207705940729
$ws.Cells.Item(10, 1).Value = "Penalty_Taking_x"
$ws.Cells.Item(10, 2).Value = 0.01252964211970782
$ws.Cells.Item(10, 3).Value = -16.59701423989972
$ws.Cells.Item(10, 4).Value = 0.01252964211970782
$ws.Cells.Item(11, 1).Value = "Vision_x"
$ws.Cells.Item(11, 2).Value = -0.1809495666157287
$ws.Cells.Item(11, 3).Value = -16.59701423989972
$ws.Cells.Item(11, 4).Value = 0.1809495666157287
$ws.Cells.Item(12, 1).Value = "Positioning_x"
$ws.Cells.Item(12, 2).Value = -0.03290121837836091
$ws.Cells.Item(12, 3).Value = -16.59701423989972
$ws.Cells.Item(12, 4).Value = 0.03290121837836091
$ws.Cells.Item(13, 1).Value = "Acceleration_x"
$ws.Cells.Item(13, 2).Value = 0.07208681128406398
$ws.Cells.Item(13, 3).Value = -16.59701423989972
$ws.Cells.Item(13, 4).Value = 0.07208681128406398
$ws.Cells.Item(14, 1).Value = "Concentration_x"
$ws.Cells.Item(14, 2).Value = -0.04060840441440411
$ws.Cells.Item(14, 3).Value = -16.59701423989972
$ws.Cells.Item(14, 4).Value = 0.04060840441440411
$ws.Cells.Item(15, 1).Value = "Pace_x"
$ws.Cells.Item(15, 2).Value = 0.08911375318663163
$ws.Cells.Item(15, 3).Value = -16.59701423989972
$ws.Cells.Item(15, 4).Value = 0.08911375318663163
$ws.Cells.Item(16, 1).Value = "Stamina_x"
$ws.Cells.Item(16, 2).Value = 0.184177793245315
$ws.Cells.Item(16, 3).Value = -16.59701423989972
$ws.Cells.Item(16, 4).Value = 0.184177793245315
$ws.Cells.Item(17, 1).Value = "Balance_x"
$ws.Cells.Item(17, 2).Value = 0.05727095918597623
$ws.Cells.Item(17, 3).Value = -16.59701423989972
$ws.Cells.Item(17, 4).Value = 0.05727095918597623
$ws.Cells.Item(18, 1).Value = "Flair_x"
$ws.Cells.Item(18, 2).Value = -0.09380265541293856
$ws.Cells.Item(18, 3).Value = -16.59701423989972
$ws.Cells.Item(18, 4).Value = 0.09380265541293856
$ws.Cells.Item(19, 1).Value = "Free_Kick_Taking_x"
$ws.Cells.Item(19, 2).Value = -0.009089981759475611
$ws.Cells.Item(19, 3).Value = -16.59701423989972
$ws.Cells.Item(19, 4).Value = 0.009089981759475611
$ws.Cells.Item(20, 1).Value = "Crossing_x"
$ws.Cells.Item(20, 2).Value = 0.129305608390409
$ws.Cells.Item(20, 3).Value = -16.59701423989972
$ws.Cells.Item(20, 4).Value = 0.129305608390409
$ws.Cells.Item(21, 1).Value = "Jumping_Reach_x"
$ws.Cells.Item(21, 2).Value = -0.05762949200016387
$ws.Cells.Item(21, 3).Value = -16.59701423989972
$ws.Cells.Item(21, 4).Value = 0.05762949200016387
$ws.Cells.Item(22, 1).Value = "Natural_Fitness_x"
$ws.Cells.Item(22, 2).Value = -0.1045181032440784
$ws.Cells.Item(22, 3).Value = -16.59701423989972
$ws.Cells.Item(22, 4).Value = 0.1045181032440784
$ws.Cells.Item(23, 1).Value = "Technique_x"
$ws.Cells.Item(23, 2).Value = 0.186089290247479
$ws.Cells.Item(23, 3).Value = -16.59701423989972
$ws.Cells.Item(23, 4).Value = 0.186089290247479
$ws.Cells.Item(24, 1).Value = "Long_Throws_x"
$ws.Cells.Item(24, 2).Value = 0.09615381609156798
$ws.Cells.Item(24, 3).Value = -16.59701423989972
$ws.Cells.Item(24, 4).Value = 0.09615381609156798
$ws.Cells.Item(25, 1).Value = "Strength_x"
$ws.Cells.Item(25, 2).Value = -0.2017506910783503
$ws.Cells.Item(25, 3).Value = -16.59701423989972
$ws.Cells.Item(25, 4).Value = 0.2017506910783503
$ws.Cells.Item(26, 1).Value = "Long_Shots_x"
$ws.Cells.Item(26, 2).Value = 0.1569422774424056
$ws.Cells.Item(26, 3).Value = -16.59701423989972
$ws.Cells.Item(26, 4).Value = 0.1569422774424056
$ws.Cells.Item(27, 1).Value = "Bravery_x"
$ws.Cells.Item(27, 2).Value = 0.05242678821507071
$ws.Cells.Item(27, 3).Value = -16.59701423989972
$ws.Cells.Item(27, 4).Value = 0.05242678821507071
$ws.Cells.Item(28, 1).Value = "Finishing_x"
$ws.Cells.Item(28, 2).Value = 0.0875766324796652
$ws.Cells.Item(28, 3).Value = -16.59701423989972
$ws.Cells.Item(28, 4).Value = 0.0875766324796652
$ws.Cells.Item(29, 1).Value = "Aggression_x"
$ws.Cells.Item(29, 2).Value = -0.09513435537874448
$ws.Cells.Item(29, 3).Value = -16.59701423989972
$ws.Cells.Item(29, 4).Value = 0.09513435537874448
$ws.Cells.Item(30, 1).Value = "Work_Rate_x"
$ws.Cells.Item(30, 2).Value = 0.0341836388601659
$ws.Cells.Item(30, 3).Value = -16.59701423989972
$ws.Cells.Item(30, 4).Value = 0.0341836388601659
$ws.Cells.Item(31, 1).Value = "Heading_x"
$ws.Cells.Item(31, 2).Value = 0.1155589284683455
$ws.Cells.Item(31, 3).Value = -16.59701423989972
$ws.Cells.Item(31, 4).Value = 0.1155589284683455
$ws.Cells.Item(32, 1).Value = "Decisions_x"
$ws.Cells.Item(32, 2).Value = -0.07593378896341374
$ws.Cells.Item(32, 3).Value = -16.59701423989972
$ws.Cells.Item(32, 4).Value = 0.07593378896341374
$ws.Cells.Item(33, 1).Value = "Tackling_x"
$ws.Cells.Item(33, 2).Value = 0.05956084540852848
$ws.Cells.Item(33, 3).Value = -16.59701423989972
$ws.Cells.Item(33, 4).Value = 0.05956084540852848
$ws.Cells.Item(34, 1).Value = "Marking_x"
$ws.Cells.Item(34, 2).Value = -0.2329563806298866
$ws.Cells.Item(34, 3).Value = -16.59701423989972
$ws.Cells.Item(34, 4).Value = 0.2329563806298866
$ws.Cells.Item(35, 1).Value = "Leadership_x"
$ws.Cells.Item(35, 2).Value = 0.02938124037041185
$ws.Cells.Item(35, 3).Value = -16.59701423989972
$ws.Cells.Item(35, 4).Value = 0.02938124037041185
$ws.Cells.Item(36, 1).Value = "Composure_x"
$ws.Cells.Item(36, 2).Value = -0.08052109668232832
$ws.Cells.Item(36, 3).Value = -16.59701423989972
$ws.Cells.Item(36, 4).Value = 0.08052109668232832
$ws.Cells.Item(37, 1).Value = "Determination_x"
$ws.Cells.Item(37, 2).Value = -0.0830455389444307
$ws.Cells.Item(37, 3).Value = -16.59701423989972
$ws.Cells.Item(37, 4).Value = 0.0830455389444307
$ws.Cells.Item(38, 1).Value = "Passing_x"
$ws.Cells.Item(38, 2).Value = -0.1423590412775249
$ws.Cells.Item(38, 3).Value = -16.59701423989972
$ws.Cells.Item(38, 4).Value = 0.1423590412775249
$ws.Cells.Item(39, 1).Value = "First_Touch_x"
$ws.Cells.Item(39, 2).Value = 0.1059821678746613
$ws.Cells.Item(39, 3).Value = -16.59701423989972
$ws.Cells.Item(39, 4).Value = 0.1059821678746613
$ws.Cells.Item(40, 1).Value = "tpr_x"
$ws.Cells.Item(40, 2).Value = 0.01272532166281879
$ws.Cells.Item(40, 3).Value = -16.59701423989972
$ws.Cells.Item(40, 4).Value = 0.01272532166281879